$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell value updates per diff ---
$ws.Range("A5").Value = 37347.45833333334
$ws.Range("A17").Value = 37712.45833333334
$ws.Range("A29").Value = 38078.45833333334
$ws.Range("A41").Value = 38443.45833333334
$ws.Range("A59").Value = 38991.45833333334
$ws.Range("C205").Value = 4373594970000
$ws.Range("D205").Value = 4373594970000
$ws.Range("E205").Value = 4373594970000
$ws.Range("F205").Value = 4373594970000
$ws.Range("C206").Value = 4411934620000
$ws.Range("D206").Value = 4411934620000
$ws.Range("E206").Value = 4411934620000
$ws.Range("F206").Value = 4411934620000
$ws.Range("C210").Value = 4577407590000
$ws.Range("D210").Value = 4577407590000
$ws.Range("E210").Value = 4577407590000
$ws.Range("F210").Value = 4577407590000
$ws.Range("C211").Value = 4566459490000
$ws.Range("D211").Value = 4566459490000
$ws.Range("E211").Value = 4566459490000
$ws.Range("F211").Value = 4566459490000
$ws.Range("C212").Value = 4592275590000
$ws.Range("D212").Value = 4592275590000
$ws.Range("E212").Value = 4592275590000
$ws.Range("F212").Value = 4592275590000
$ws.Range("C213").Value = 4639859400000
$ws.Range("D213").Value = 4639859400000
$ws.Range("E213").Value = 4639859400000
$ws.Range("F213").Value = 4639859400000
$ws.Range("C214").Value = 4641345140000
$ws.Range("D214").Value = 4641345140000
$ws.Range("E214").Value = 4641345140000
$ws.Range("F214").Value = 4641345140000
$ws.Range("C215").Value = 4681223420000
$ws.Range("D215").Value = 4681223420000
$ws.Range("E215").Value = 4681223420000
$ws.Range("F215").Value = 4681223420000
$ws.Range("C216").Value = 4725508480000
$ws.Range("D216").Value = 4725508480000
$ws.Range("E216").Value = 4725508480000
$ws.Range("F216").Value = 4725508480000
$ws.Range("C217").Value = 4680322510000
$ws.Range("D217").Value = 4680322510000
$ws.Range("E217").Value = 4680322510000
$ws.Range("F217").Value = 4680322510000
$ws.Range("C219").Value = 4809150480000
$ws.Range("D219").Value = 4809150480000
$ws.Range("E219").Value = 4809150480000
$ws.Range("F219").Value = 4809150480000
$ws.Range("C221").Value = 4958595660000
$ws.Range("D221").Value = 4958595660000
$ws.Range("E221").Value = 4958595660000
$ws.Range("F221").Value = 4958595660000
$ws.Range("C222").Value = 5004666910000
$ws.Range("D222").Value = 5004666910000
$ws.Range("E222").Value = 5004666910000
$ws.Range("F222").Value = 5004666910000
$ws.Range("C223").Value = 5020790900000
$ws.Range("D223").Value = 5020790900000
$ws.Range("E223").Value = 5020790900000
$ws.Range("F223").Value = 5020790900000
$ws.Range("C224").Value = 5059232680000
$ws.Range("D224").Value = 5059232680000
$ws.Range("E224").Value = 5059232680000
$ws.Range("F224").Value = 5059232680000
$ws.Range("C225").Value = 5094308060000
$ws.Range("D225").Value = 5094308060000
$ws.Range("E225").Value = 5094308060000
$ws.Range("F225").Value = 5094308060000
$ws.Range("C226").Value = 5178041490000
$ws.Range("D226").Value = 5178041490000
$ws.Range("E226").Value = 5178041490000
$ws.Range("F226").Value = 5178041490000
$ws.Range("C227").Value = 5214187690000
$ws.Range("D227").Value = 5214187690000
$ws.Range("E227").Value = 5214187690000
$ws.Range("F227").Value = 5214187690000
$ws.Range("C228").Value = 5235568230000
$ws.Range("D228").Value = 5235568230000
$ws.Range("E228").Value = 5235568230000
$ws.Range("F228").Value = 5235568230000
$ws.Range("C229").Value = 5179738620000
$ws.Range("D229").Value = 5179738620000
$ws.Range("E229").Value = 5179738620000
$ws.Range("F229").Value = 5179738620000
$ws.Range("C230").Value = 5290478980000
$ws.Range("D230").Value = 5290478980000
$ws.Range("E230").Value = 5290478980000
$ws.Range("F230").Value = 5290478980000
$ws.Range("C231").Value = 5390398340000
$ws.Range("D231").Value = 5390398340000
$ws.Range("E231").Value = 5390398340000
$ws.Range("F231").Value = 5390398340000
$ws.Range("C232").Value = 5449356120000
$ws.Range("D232").Value = 5449356120000
$ws.Range("E232").Value = 5449356120000
$ws.Range("F232").Value = 5449356120000
$ws.Range("C233").Value = 5471474170000
$ws.Range("D233").Value = 5471474170000
$ws.Range("E233").Value = 5471474170000
$ws.Range("F233").Value = 5471474170000
$ws.Range("C234").Value = 5507491430000
$ws.Range("D234").Value = 5507491430000
$ws.Range("E234").Value = 5507491430000
$ws.Range("F234").Value = 5507491430000
$ws.Range("C236").Value = 5564521500000
$ws.Range("D236").Value = 5564521500000
$ws.Range("E236").Value = 5564521500000
$ws.Range("F236").Value = 5564521500000
$ws.Range("C239").Value = 5617130550000
$ws.Range("D239").Value = 5617130550000
$ws.Range("E239").Value = 5617130550000
$ws.Range("F239").Value = 5617130550000
$ws.Range("C240").Value = 5647837280000
$ws.Range("D240").Value = 5647837280000
$ws.Range("E240").Value = 5647837280000
$ws.Range("F240").Value = 5647837280000
$ws.Range("C241").Value = 5542014840000
$ws.Range("D241").Value = 5542014840000
$ws.Range("E241").Value = 5542014840000
$ws.Range("F241").Value = 5542014840000
$ws.Range("C242").Value = 5630383690000
$ws.Range("D242").Value = 5630383690000
$ws.Range("E242").Value = 5630383690000
$ws.Range("F242").Value = 5630383690000
$ws.Range("C243").Value = 5704249840000
$ws.Range("D243").Value = 5704249840000
$ws.Range("E243").Value = 5704249840000
$ws.Range("F243").Value = 5704249840000
$ws.Range("C244").Value = 5739159050000
$ws.Range("D244").Value = 5739159050000
$ws.Range("E244").Value = 5739159050000
$ws.Range("F244").Value = 5739159050000
$ws.Range("C245").Value = 5742427260000
$ws.Range("D245").Value = 5742427260000
$ws.Range("E245").Value = 5742427260000
$ws.Range("F245").Value = 5742427260000
$ws.Range("C246").Value = 5825723830000
$ws.Range("D246").Value = 5825723830000
$ws.Range("E246").Value = 5825723830000
$ws.Range("F246").Value = 5825723830000
$ws.Range("C247").Value = 5801917230000
$ws.Range("D247").Value = 5801917230000
$ws.Range("E247").Value = 5801917230000
$ws.Range("F247").Value = 5801917230000
$ws.Range("C248").Value = 5833040250000
$ws.Range("D248").Value = 5833040250000
$ws.Range("E248").Value = 5833040250000
$ws.Range("F248").Value = 5833040250000
$ws.Range("C249").Value = 5855415460000
$ws.Range("D249").Value = 5855415460000
$ws.Range("E249").Value = 5855415460000
$ws.Range("F249").Value = 5855415460000
$ws.Range("C250").Value = 5887405600000
$ws.Range("D250").Value = 5887405600000
$ws.Range("E250").Value = 5887405600000
$ws.Range("F250").Value = 5887405600000
$ws.Range("C251").Value = 5915934540000
$ws.Range("D251").Value = 5915934540000
$ws.Range("E251").Value = 5915934540000
$ws.Range("F251").Value = 5915934540000
$ws.Range("C252").Value = 5940210650000
$ws.Range("D252").Value = 5940210650000
$ws.Range("E252").Value = 5940210650000
$ws.Range("F252").Value = 5940210650000
$ws.Range("C256").Value = 6077524080000
$ws.Range("D256").Value = 6077524080000
$ws.Range("E256").Value = 6077524080000
$ws.Range("F256").Value = 6077524080000
$ws.Range("C258").Value = 6224248910000
$ws.Range("D258").Value = 6224248910000
$ws.Range("E258").Value = 6224248910000
$ws.Range("F258").Value = 6224248910000

# --- Append new row 259 (new data point) ---
# Copy formatting from row 258 first so the new row matches existing styling
$ws.Range("A258:G258").Copy()
$ws.Range("A259:G259").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A259").Value = 45078.41666666666
$ws.Range("B259").Value = "ECONOMICS:CZM2"
$ws.Range("C259").Value = 6243183470000
$ws.Range("D259").Value = 6243183470000
$ws.Range("E259").Value = 6243183470000
$ws.Range("F259").Value = 6243183470000
$ws.Range("G259").Value = 0
